# Updated symbol list on Mon Jan  9 09:27:18 UTC 2023 with GitHub Actions
#
# The sheet is a scraped "coin ranking" table. Columns:
#   A=index(#) B=Coin C=Link D=Price E=Volume(1h) F=Data G=Hora
# This refresh: (1) re-pulled Price/Volume for most rows, and
# (2) CoinExToken jumped from rank #21 (row 23) up to rank #13 (row 15),
# pushing One/TigerCash/LEO/GateToken/BTSEToken/BitpandaEcosystemToken/
# ProBitToken/MCDex each down one row (their Coin+Link shift, Price/Volume
# are refreshed to the newly-pulled figures for that rank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as literal text, even when the text
# looks numeric (e.g. "278.28" or "6.51%") so Excel doesn't silently
# coerce it into a Number/Percentage cell. We briefly force the Text
# number format, assign the literal, then restore the "Normal" style so
# no stray formatting is left behind on the cell.
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row -> column -> new value. Only the cells actually touched by this
# refresh are listed; B/C are plain (non-numeric-looking) text so they're
# set directly, D/E go through Set-TextValue since they hold numeric- or
# percent-looking strings that must stay text.
$rowUpdates = [ordered]@{
    2  = @{ D = '278.28';   E = '6.51%' }
    3  = @{ D = '27.31';    E = '1.45%' }
    4  = @{ D = '4.789';    E = '1.50%' }
    5  = @{ D = '0.06287';  E = '1.14%' }
    6  = @{ D = '6.904';    E = '2.61%' }
    7  = @{ D = '0.8767';   E = '3.22%' }
    8  = @{ D = '0.9443';   E = '3.54%' }
    9  = @{ D = '0.1458';   E = '4.08%' }
    10 = @{ D = '0.05170';  E = '4.21%' }
    11 = @{ D = '0.07280';  E = '2.71%' }
    12 = @{ D = '0.03103';  E = '-0.34%' }
    13 = @{ D = '0.09071';  E = '0.17%' }
    14 = @{ D = '0.001563'; E = '2.05%' }
    15 = @{ B = 'CoinExToken';            C = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet';            D = '0.04330';     E = '2.43%' }
    16 = @{ B = 'One';                    C = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one';                     D = '0.0006272';   E = '1.82%' }
    17 = @{ B = 'TigerCash';              C = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch';                   D = '0.005877';    E = '-1.67%' }
    18 = @{ B = 'LEO';                    C = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo';                     D = '3.446';       E = '-0.07%' }
    19 = @{ B = 'GateToken';              C = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt';                D = '3.266';       E = '2.96%' }
    20 = @{ B = 'BTSEToken';              C = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse';              D = '2.285';       E = '5.39%' }
    21 = @{ B = 'BitpandaEcosystemToken'; C = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best';     D = '0.3147';      E = '1.61%' }
    22 = @{ B = 'ProBitToken';            C = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob';                D = '0.1312';      E = '0.13%' }
    23 = @{ B = 'MCDex';                  C = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb';                      D = '3.855';       E = '-6.25%' }
    24 = @{ D = '0.001180';  E = '-0.18%' }
    25 = @{ D = '0.004279';  E = '5.10%' }
    26 = @{ D = '0.0001200'; E = '-0.03%' }
    27 = @{ E = '3.16%' }
    40 = @{ D = '0.04069';     E = '3.15%' }
    41 = @{ D = '0.006530';    E = '57.93%' }
    42 = @{ D = '0.1157';      E = '3.99%' }
    43 = @{ D = '0.002175';    E = '1.54%' }
    44 = @{ D = '0.01185';     E = '-9.90%' }
    45 = @{ D = '0.00005130';  E = '-0.64%' }
    46 = @{ E = '0.03%' }
    47 = @{ E = '856.95%' }
    48 = @{ D = '0.02252';     E = '-33.81%' }
    49 = @{ D = '0.00002102';  E = '0.03%' }
    50 = @{ E = '0.03%' }
}

foreach ($row in $rowUpdates.Keys) {
    $cols = $rowUpdates[$row]
    if ($cols.ContainsKey('B')) { $ws.Range("B$row").Value = $cols['B'] }
    if ($cols.ContainsKey('C')) { $ws.Range("C$row").Value = $cols['C'] }
    if ($cols.ContainsKey('D')) { Set-TextValue "D$row" $cols['D'] }
    if ($cols.ContainsKey('E')) { Set-TextValue "E$row" $cols['E'] }
}
